# The workbook is already open; grab the workbook/active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the two test-step actions in row 2 (Action1 / Action2 columns) ---
# "login"  -> "productCatalogPage"
# "logout" -> "applyFilters"
$ws.Range("D2").Value = "productCatalogPage"
$ws.Range("E2").Value = "applyFilters"

# --- Resize the data columns (A:D) and give column E an explicit width too,
#     matching the layout recorded after the edit ---
$ws.Columns.Item(1).ColumnWidth = 18.8333333333333
$ws.Columns.Item(2).ColumnWidth = 16.6666666666667
$ws.Columns.Item(3).ColumnWidth = 17.0
$ws.Columns.Item(4).ColumnWidth = 24.5
$ws.Columns.Item(5).ColumnWidth = 7.5

# --- Move the active selection from C2 to C3 ---
$ws.Range("C3").Select()
